$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 11466
$ws.Range("C3").Value = 6147
$ws.Range("C4").Value = 9155
$ws.Range("C5").Value = 6886
$ws.Range("C6").Value = 5270
$ws.Range("C7").Value = 8567
$ws.Range("C8").Value = 22879
$ws.Range("C9").Value = 16927
$ws.Range("C10").Value = 5087
$ws.Range("C11").Value = 3989
$ws.Range("C12").Value = 57
